$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1:W1").AutoFilter()

$sortRange = $ws.Range("A1:W19")
$key1 = $ws.Range("A1")

$sortRange.Sort($key1, 1, $null, $null, 2, $null, 1, 1)

$ws.Range("W12").Select()
